$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(155).Insert()

$ws.Cells.Item(155, 1).Value = 4
$ws.Cells.Item(155, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(155, 3).Value = "Los Lagos"
$ws.Cells.Item(155, 4).Value = 44508
$ws.Cells.Item(155, 5).Value = 10
$ws.Cells.Item(155, 6).Value = 100112045
$ws.Cells.Item(155, 7).Value = "Zapallo"
$ws.Cells.Item(155, 8).Value = "Paine"
$ws.Cells.Item(155, 9).Value = "1a (guarda)"
$ws.Cells.Item(155, 10).Value = 500
$ws.Cells.Item(155, 11).Value = 400
$ws.Cells.Item(155, 12).Value = 400
$ws.Cells.Item(155, 13).Value = 400
$ws.Cells.Item(155, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(155, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(155, 16).Value = 400
$ws.Cells.Item(155, 17).Value = 1
$ws.Cells.Item(155, 18).Value = "Hortaliza"
